# Generate Report for Handoff
# Refresh the localization-status report for the new handoff cycle:
# new source-file GUID, new xliff hashes, and updated handoff/xliff
# generation timestamps.

$wb = $excel.ActiveWorkbook

$oldGuid = "480c67af-41a9-417c-9149-0ed1580777a4"
$newGuid = "49cc0501-330e-4578-af8d-66c66e9a4998"

$oldHash = "2736ff5ea29c4b47f35a5b32f852da646bb33326"
$newHash = "2309ce991a722ca8273d5d0894b5ceb849f124a9"

# The hyperlink's target URL is unchanged across this edit - only the
# displayed text (file name) needs to be refreshed to the new GUID.
$hyperlinkUrl = "https://github.com/OpenLocalizationTestOrg/oltest/blob/6c7d81003e4877bd809d887590a173b71eefe05b/e2e/$oldGuid.md"

# --- Sheet "Overview" ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("B2").Value = "e2e\$newGuid.md"
$rngOverviewB2 = $wsOverview.Range("B2")
$rngOverviewB2.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($rngOverviewB2, $hyperlinkUrl, "", "", "e2e\$newGuid.md") | Out-Null
$wsOverview.Range("G2").Value = "2016-08-12 09:09:54"

# --- Sheet "zh-cn" ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("A2").Value = "$newGuid.md"
$rngZhA2 = $wsZh.Range("A2")
$rngZhA2.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($rngZhA2, $hyperlinkUrl, "", "", "$newGuid.md") | Out-Null
$wsZh.Range("G2").Value = "$newGuid.$newHash.zh-cn.xlf"
$wsZh.Range("H2").Value = "2016-08-12 09:09:47"

# --- Sheet "de-de" ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("A2").Value = "$newGuid.md"
$rngDeA2 = $wsDe.Range("A2")
$rngDeA2.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($rngDeA2, $hyperlinkUrl, "", "", "$newGuid.md") | Out-Null
$wsDe.Range("G2").Value = "$newGuid.$newHash.de-de.xlf"
$wsDe.Range("H2").Value = "2016-08-12 09:09:54"
